$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for data rows 2-42
# from serial date 45716 (2025-02-28) to 45717 (2025-03-01)
for ($row = 2; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45716) {
        $cell.Value2 = 45717
    }
}
